$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the two "live" credential rows with the new username/password pair.
$ws.Range("A2").Value = "lennox.fernandes@gmail.com"
$ws.Range("B2").Value = "abc@123"

# The remaining sample rows (3-6) are wiped out in the new sheet.
$ws.Range("A3:B6").ClearContents()

# Turn the row-2 credentials into mailto hyperlinks (this also applies the
# built-in "Hyperlink" cell style to A2/B2, matching the new styles.xml).
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:lennox.fernandes@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:abc@123")

# Selection moved to A8 before the file was saved.
$ws.Range("A8").Select() | Out-Null
